$d = $word.ActiveDocument

# Remove the stray "$" character that was accidentally left after the
# "Recorder: Records all answers and questions." bullet point.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("questions.$", $true, $false, $false, $false, $false, `
               $true, 1, $false, "questions.", 2)
